# Rewrite speaker notes to concise reference style (craft-skill convention).
# Touches the "Notes Placeholder" body shape on the Notes Page of 5 slides
# in this deck; all other content is left untouched.

$p = $ppt.ActivePresentation

function Set-NotesText {
    param($Presentation, [int]$SlideIndex, [string]$NewText)

    $slide = $Presentation.Slides.Item($SlideIndex)
    $notesPage = $slide.NotesPage

    for ($i = 1; $i -le $notesPage.Shapes.Count; $i++) {
        $shape = $notesPage.Shapes.Item($i)
        if ($shape.Name -like "Notes Placeholder*") {
            $shape.TextFrame.TextRange.Text = $NewText
            return
        }
    }
}

# Slide 5 — contrast stat callout (source line is unchanged; only the
# second line, after the embedded line break, is rewritten).
Set-NotesText $p 5 "Source: Addy Osmani — https://addyo.substack.com/p/the-reality-of-ai-assisted-software`nKey contrast: perceived 20% faster vs actual 19% slower after debugging. Gap = active harm, not just missed potential. ~15 seconds."

# Slide 6 — "let it land" beat.
Set-NotesText $p 6 "Let it land. Goal: self-identification. ~10 seconds."

# Slide 7 — spectrum self-identification beat.
Set-NotesText $p 7 "Self-identification moment. Pause for audience to locate themselves on spectrum. ~20 seconds."

# Slide 9 — session walkthrough beat.
Set-NotesText $p 9 "Brief walkthrough of each session. One-liners = curiosity hooks. ~30 seconds."

# Slide 10 — low-commitment ask beat.
Set-NotesText $p 10 "Specific, low-commitment ask. Emphasis: existing weekly task, not new skill. Live transformation in session. ~15 seconds."
